# The deck orders slides as: ... slide7 ("The New" / Routing / Network Topology),
# slide8 ("The New" / Problems, with the OpenFlow sentence) ...
$p = $ppt.ActivePresentation

# --- Slide 7: add a new "Define SDN" paragraph before "Routing" ---
$slide7 = $p.Slides.Item(7)
$contentShape7 = $slide7.Shapes.Item(1)
$tr7 = $contentShape7.TextFrame.TextRange
$tr7.Text = "Define SDN`rRouting`rNetwork Topology"

# --- Slide 8: change the third run of the OpenFlow paragraph ---
$slide8 = $p.Slides.Item(8)
$contentShape8 = $slide8.Shapes.Item(1)
$tr8 = $contentShape8.TextFrame.TextRange
$openFlowPara = $tr8.Paragraphs(2)
$lastRun = $openFlowPara.Runs(3)
$lastRun.Text = ", Define it)"
